$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.284.92'
$ws.Range('D3').Value = '1.859.41'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''0.7039'
$ws.Range('E5').Value = '  +1.07%  '
$ws.Range('D6').Value = '''238.21'
$ws.Range('E6').Value = '  +0.56%  '
$ws.Range('D7').Value = '''1.000'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '''0.07924'
$ws.Range('E8').Value = '  +2.88%  '
$ws.Range('D9').Value = '''0.3043'
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('D10').Value = '''24.60'
$ws.Range('E10').Value = '  +6.05%  '
$ws.Range('D11').Value = '''0.08182'
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('D12').Value = '''5.222'
$ws.Range('E12').Value = '  +1.55%  '
$ws.Range('D13').Value = '''0.7170'
$ws.Range('E13').Value = '  +0.29%  '
$ws.Range('E14').Value = '  +0.65%  '
$ws.Range('D15').Value = '1.696.02'
$ws.Range('E15').Value = '  -8.46%  '
$ws.Range('D16').Value = '28.839.32'
$ws.Range('E16').Value = '  -1.23%  '
$ws.Range('D17').Value = '''5.817'
$ws.Range('E17').Value = '  +1.40%  '
$ws.Range('D18').Value = '''0.000007797'
$ws.Range('E18').Value = '  +0.89%  '
$ws.Range('D19').Value = '''13.22'
$ws.Range('E19').Value = '  -0.51%  '
$ws.Range('D20').Value = '''238.22'
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('D21').Value = '''0.9998'
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').Value = '''1.000'
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').Value = '''7.541'
$ws.Range('E23').Value = '  +1.74%  '
$ws.Range('D24').Value = '1.930.01'
$ws.Range('E24').Value = '  -7.84%  '
$ws.Range('D25').Value = '''162.51'
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').Value = '''8.892'
$ws.Range('E26').Value = '  -1.07%  '
$ws.Range('D27').Value = '''0.1426'
$ws.Range('E27').Value = '  -3.85%  '
$ws.Range('D28').Value = '''18.09'
$ws.Range('D29').Value = '''1.918'
$ws.Range('E29').Value = '  -5.90%  '
$ws.Range('D30').Value = '''1.378'
$ws.Range('E30').Value = '  -2.75%  '
$ws.Range('D31').Value = '''1.477'
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('D32').Value = '''4.325'
$ws.Range('E32').Value = '  -2.18%  '
$ws.Range('D33').Value = '''4.060'
$ws.Range('E33').Value = '  +1.25%  '
$ws.Range('D34').Value = '''0.05174'
$ws.Range('E34').Value = '  -0.21%  '
$ws.Range('D35').Value = '''1.177'
$ws.Range('E35').Value = '  +1.35%  '
$ws.Range('D36').Value = '''0.7105'
$ws.Range('E36').Value = '  +0.37%  '
$ws.Range('D37').Value = '''0.9901'
$ws.Range('E37').Value = '  -0.97%  '
$ws.Range('D38').Value = '''2.675'
$ws.Range('E38').Value = '  +0.73%  '
$ws.Range('D39').Value = '''0.01851'
$ws.Range('E39').Value = '  +0.47%  '
$ws.Range('D40').Value = '''2.687'
$ws.Range('E40').Value = '  -1.42%  '
$ws.Range('D41').Value = '1.156.70'
$ws.Range('E41').Value = '  +1.42%  '
$ws.Range('D42').Value = '''0.9246'
$ws.Range('E42').Value = '  -1.62%  '
$ws.Range('D43').Value = '''5.940'
$ws.Range('E43').Value = '  +1.11%  '
$ws.Range('D44').Value = '''0.4252'
$ws.Range('E44').Value = '  -0.30%  '
$ws.Range('E45').Value = '  -0.94%  '
$ws.Range('D46').Value = '''0.9996'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').Value = '''100.85'
$ws.Range('E47').Value = '  -2.22%  '
$ws.Range('D48').Value = '''0.5315'
$ws.Range('E48').Value = '  -2.83%  '
$ws.Range('D49').Value = '''1.755'
$ws.Range('E49').Value = '  -2.00%  '
$ws.Range('D50').Value = '''9.178'
$ws.Range('E50').Value = '  +0.35%  '
$ws.Range('D51').Value = '''6.997'
$ws.Range('E51').Value = '  +1.00%  '

# Reset style on the Price column range to clear any quote-prefix / text-format
# artifacts introduced by assigning numeric-looking strings as text, restoring
# the original default (General) appearance for all cells in that range.
$ws.Range("D2:D51").Style = "Normal"
